$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Remove the stale "_GoBack" bookmark that currently sits near the top
#    of the document (in the paragraph that also holds the floating
#    text-box drawing with the "Виконав/Перевірив" signature block).
#    Word always keeps exactly one "_GoBack" bookmark, tracking the
#    location of the most recent edit; since the edit below happens in
#    the conclusion paragraph, the bookmark must move there. We delete
#    the old one now and add the new one after editing that paragraph.
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ---------------------------------------------------------------------
# 2. Rewrite the conclusion paragraph's text: shorten/re-word it.
# ---------------------------------------------------------------------
$oldText = "Завдяки цій лабораторній роботі я дізнався про декілька типів зв’язаних списків (лінійний, двозв’язний,) навчився їх створювати, видаляти та використовувати. Також я більш досконало вивчив використання вказівників, їх запис та обробка в структурах. Під час виконання лабораторної роботи я використовував додаткові файли для зчитування посимвольно, що також потребувало додаткової інформації. Найбільший інтерес викликали ітеративні цикли з використанням вказівників, а не змінних. Ця лабораторна робота, на відміну від попередніх, була виконана найбільш стисло, оскільки включає у себе лише три функції для оперування списком."
$newText = "Завдяки цій лабораторній роботі я дізнався про такий тип списоку як однозв’язний навчився його створювати, видаляти та використовувати. Також я більш досконало вивчив використання вказівників, їх запис та обробка в структурах. Під час виконання лабораторної роботи я використовував додаткові файли для зчитування посимвольно, що також потребувало додаткової інформації."

$found = $d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
if (-not $found) {
    Write-Host "WARNING: conclusion paragraph text was not found / replaced."
}

# ---------------------------------------------------------------------
# 3. The conclusion paragraph used to be followed by a run that contains
#    only a single trailing space (" "); that run is the very last piece
#    of content in the document's main story, right before the paragraph
#    mark. Select it precisely by character offset, wrap it with the new
#    "_GoBack" bookmark, then delete it. Deleting the (now bookmarked)
#    space collapses the bookmark to an empty range that naturally lands
#    right after the rewritten run -- exactly where Word itself leaves
#    "_GoBack" after the last edit.
# ---------------------------------------------------------------------
$bodyEnd = $d.Content.End
$pilcrowStart = $bodyEnd - 1
$spaceStart = $pilcrowStart - 1
$spaceRange = $d.Range($spaceStart, $pilcrowStart)

if ($spaceRange.Text -eq " ") {
    $d.Bookmarks.Add("_GoBack", $spaceRange)
    $spaceRange.Delete()
} else {
    Write-Host "WARNING: trailing space run not where expected, got:" $spaceRange.Text
    # Fall back: just drop a fresh bookmark at the end of the main story.
    $endRange = $d.Range($d.Content.End - 1, $d.Content.End - 1)
    $d.Bookmarks.Add("_GoBack", $endRange)
}
